# Update column F (dSF) values on specific rows to match repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    10 = -2
    15 = -1
    17 = -2
    18 = 1
    23 = 2
    24 = -1
    26 = -1
    34 = -2
    35 = 1
    41 = -2
    43 = 2
    56 = -9
    61 = 3
    65 = -1
    69 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
